# Update the "想去人数" (interest count) values in the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 1462
$wsExpo.Range("F3").Value = 3074
$wsExpo.Range("F4").Value = 40
$wsExpo.Range("F5").Value = 615
$wsExpo.Range("F6").Value = 289

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 1462
$wsAll.Range("F3").Value = 3074
$wsAll.Range("F4").Value = 40
$wsAll.Range("F5").Value = 615
$wsAll.Range("F7").Value = 289
